$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2021" data column (R) entirely - it was the last column
# (A:R) on the sheet, so deleting it shrinks the used range to A:Q and
# drops every R-column cell (R4:R14) along with their values.
$ws.Range("R:R").Delete() | Out-Null

# The author's selection moved to N19 after the edit.
$ws.Range("N19").Select() | Out-Null
